# "add auxiliary functions and storages"
#
# Sheet1 gains a third data series ("0-zaxis") in column C. The values
# that used to sit in column C (for the rows that already had a partial
# third column) move over to become the new column B values, and column
# C is filled down with a constant 40 for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New header for column C, matching the header style already used by
# A1/B1 ("0-xaxis"/"0-yaxis").
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("C1").Value = "0-zaxis"

# The handful of column-B values that were previously mirrored (in a
# shifted form) over in column C now take on those values directly.
$ws.Range("B3").Value = 10.208
$ws.Range("B4").Value = 15.292
$ws.Range("B5").Value = 20.457999999999998
$ws.Range("B11").Value = 38.167000000000002

# Column C becomes a uniform auxiliary storage column, value 40 for
# every data row (2-11).
for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 3).Value = 40
}

# Mirror the saved cursor position recorded for this edit.
$ws.Range("E6").Select()
